# eShopping-user-stories.xlsx -- "Add files via upload" edit
#
# Semantic changes (once the shared-string re-indexing in the diff is
# resolved against actual text) are three feature rows on Sheet1:
#   - ESUC003 (row 6):  Feature Owner  "Scrum Master" -> "Mohd.Areeb"
#   - ESUC004 (row 7):  Short Description, Dependency and Feature Owner
#                        replaced with new "Describing the summary..." story
#   - ESUC010 (row 13): Short Description, Dependency and Feature Owner
#                        replaced with new "Register user with the
#                        credentials" story, and the row's custom (wrapped)
#                        height is no longer needed so it reverts to default.
# The workbook's view also ends up scrolled to a different cell/column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 (ESUC003): only the Feature Owner changes -----------------------
$ws.Range("F6").Value = "Mohd.Areeb"

# --- Row 7 (ESUC004): new user story text -----------------------------------
$ws.Range("D7").Value = "Describing the summary of the products"
$ws.Range("E7").Value = "I can see the elaborate description of the products"
$ws.Range("F7").Value = "Pushkar sinha"

# --- Row 13 (ESUC010): new user story text ----------------------------------
$ws.Range("D13").Value = "Register user with the credentials"
$ws.Range("E13").Value = "Access the full functionality of the webapp"
$ws.Range("F13").Value = "Aryaman garg"

# The old text needed two lines (ht="26"); the replacement text fits on one
# line, so Excel drops the explicit row height back to the sheet default.
$ws.Rows.Item(13).AutoFit()

# --- View state: re-freeze header row and move the selection ----------------
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("F13").Select()
